# Add a new vendor row (Netflix) to the Vendors sheet, matching the
# existing pattern used for the other vendor rows (Amazon, Uber Eats,
# carvana, Best Buy, Walmart).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new row. B7 is written before A7 so the shared-string
# table picks up "www.netflix.com" ahead of "Netflix", matching the
# order vendor web addresses were entered for the other rows.
$ws.Range("B7").Value = "www.netflix.com"
$ws.Range("A7").Value = "Netflix"

# Give the new web-address cell the same "Hyperlink" cell style used by
# the other vendor links, and create the actual hyperlink relationship.
$ws.Range("B7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B7"), "http://www.netflix.com/")
$ws.Range("B7").Style = "Hyperlink"

# Move the active selection the way it ended up after the edit.
[void]$ws.Range("A9").Select()
